$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 18:02:51"
$wsZhCn.Range("G5").Value = "2016-02-22 18:03:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 18:03:04"
$wsDeDe.Range("G5").Value = "2016-02-22 18:03:54"
